$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.980.91"
$ws.Range("E2").Value = "  +1.17%  "
# Row 3
$ws.Range("D3").Value = "1.895.14"
$ws.Range("E3").Value = "  +0.68%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.017"
$ws.Range("E4").Value = "  +1.59%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.87"
$ws.Range("E5").Value = "  +1.46%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.016"
$ws.Range("E6").Value = "  +1.45%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4700"
$ws.Range("E7").Value = "  -0.28%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3938"
$ws.Range("E8").Value = "  -0.65%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.70"
$ws.Range("E9").Value = "  -1.59%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08074"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.022"
$ws.Range("E11").Value = "  -0.44%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.87"
$ws.Range("E12").Value = "  +0.25%  "
# Row 13
$ws.Range("D13").Value = "1.894.11"
$ws.Range("E13").Value = "  +0.62%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.981"
$ws.Range("E14").Value = "  +0.30%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.154"
# Row 16
$ws.Range("E16").Value = "  +1.51%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06785"
$ws.Range("E17").Value = "  +2.81%  "
# Row 18
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.76"
$ws.Range("E18").Value = "  +1.00%  "
# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001055"
$ws.Range("E19").Value = "  +1.09%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.24"
$ws.Range("E20").Value = "  +0.23%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.015"
$ws.Range("E21").Value = "  +1.36%  "
# Row 22
$ws.Range("D22").Value = "27.999.41"
$ws.Range("E22").Value = "  +1.20%  "
# Row 23
$ws.Range("E23").Value = "  +0.17%  "
# Row 24
$ws.Range("E24").Value = "  +0.22%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.346"
$ws.Range("E25").Value = "  +1.79%  "
# Row 26
$ws.Range("D26").Value = "2.112.60"
$ws.Range("E26").Value = "  +0.40%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.57"
$ws.Range("E27").Value = "  +2.97%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.10"
$ws.Range("E28").Value = "  -0.75%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.090"
$ws.Range("E29").Value = "  -0.36%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.487"
$ws.Range("E30").Value = "  -1.81%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.14"
$ws.Range("E31").Value = "  -0.41%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9742"
$ws.Range("E32").Value = "  +1.24%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09519"
$ws.Range("E33").Value = "  -0.14%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.656"
$ws.Range("E34").Value = "  +1.08%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.406"
$ws.Range("E35").Value = "  -4.77%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.382"
$ws.Range("E36").Value = "  +1.57%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06141"
$ws.Range("E37").Value = "  +0.37%  "
# Row 38
$ws.Range("E38").Value = "  +0.27%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.221"
$ws.Range("E39").Value = "  -0.31%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.123"
$ws.Range("E40").Value = "  -0.67%  "
# Row 41
$ws.Range("E41").Value = "  +0.02%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1894"
$ws.Range("E42").Value = "  -0.24%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.33"
$ws.Range("E43").Value = "  -0.08%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.270"
$ws.Range("E44").Value = "  +1.73%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5720"
$ws.Range("E45").Value = "  +0.15%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.32"
$ws.Range("E46").Value = "  +1.42%  "
# Row 47
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.408"
$ws.Range("E47").Value = "  -0.12%  "
# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.945"
$ws.Range("E48").Value = "  +0.35%  "
# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06939"
$ws.Range("E49").Value = "  +1.64%  "
# Row 50
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.09"
$ws.Range("E50").Value = "  +3.46%  "
# Row 51
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.072"
$ws.Range("E51").Value = "  +0.47%  "
